# Update Training Dashboard progress as of date 04-Nov-2025
# For rows 3-19: column H (PERIOD TO EXPIRE) decreases by 1,
# column I (LAST UPDATE) changes from 03-Nov-2025 to 04-Nov-2025.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Dashboard")

$xlPasteValues = -4163

for ($row = 3; $row -le 19; $row++) {
    $hCell = $ws.Cells.Item($row, 8)   # Column H - PERIOD TO EXPIRE
    $iCell = $ws.Cells.Item($row, 9)   # Column I - LAST UPDATE

    # Decrement the numeric "period to expire" value by 1 day.
    $hCell.Value2 = $hCell.Value2 - 1

    # Write the new date as literal text (not an auto-converted date
    # serial) by computing it via a formula and then freezing the
    # formula result back down to a plain value in place.
    $iCell.Formula = '="04-Nov-2025"'
    $iCell.Copy()
    $iCell.PasteSpecial($xlPasteValues)
}

$excel.CutCopyMode = $false
